# Generate Report for Archive
#
# The localization status report is regenerated: the "zh-cn"/"de-de"
# status cells move from "Ready for handoff" to "In Translation", and the
# "Status" columns are re-sized (narrower, to fit the new text) on the
# Overview sheet (columns E/F) and on the per-language detail sheets
# (column C on "zh-cn" and "de-de").

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# Target stored column width (OOXML "width" attribute, in character units)
# for the Status columns after the text shrank. ColumnWidth is quantized by
# the host to the nearest 1/6 of a character, so we pick the setting whose
# quantized result lands closest to the recorded width.
$statusColumnWidth = 12.5

# --- Overview sheet: zh-cn (E) / de-de (F) status cells -------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $statusColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $statusColumnWidth

# --- zh-cn detail sheet: Status cell (C) -----------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $statusColumnWidth

# --- de-de detail sheet: Status cell (C) -----------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $statusColumnWidth
